# Update scripts with new TPM (tpm) values for Lta-Tnfrsf14 LR pair.
# Rows 2-11 get refreshed numeric/category values, and 5 new rows (12-16)
# are appended for the "Resolving-Mac" sending-cluster group, mirroring the
# existing per-sending-cluster block of 5 target clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2..16, columns A..T:
#   A: Sending cluster, B: Ligand symbol, C: Receptor symbol, D: Target cluster,
#   E..T: the associated NATMI TPM-derived metrics.
$data = @(
    @("ECs","Lta","Tnfrsf14","ECs",3,1,0.6836903333333334,2.051071,0.7427603292077943,0.7427603292077943,3,1,5.578493666666667,16.735481,0.1036332930693284,0.1036332930693284,3.813962194461223,34.325659750151,0.07697469887706218,0.07697469887706221),
    @("ECs","Lta","Tnfrsf14","FAPs",3,1,0.6836903333333334,2.051071,0.7427603292077943,0.7427603292077943,3,1,3.704475666666667,11.113427,0.06881911773528272,0.06881911773528274,2.532714203368556,22.794427830317,0.05111611054484855,0.05111611054484856),
    @("ECs","Lta","Tnfrsf14","Inflammatory-Mac",3,1,0.6836903333333334,2.051071,0.7427603292077943,0.7427603292077943,3,1,24.77295966666667,74.31887900000001,0.4602144490493554,0.4602144490493556,16.93703305215656,152.433297469409,0.3418290356820829,0.3418290356820831),
    @("ECs","Lta","Tnfrsf14","MuSCs",3,1,0.6836903333333334,2.051071,0.7427603292077943,0.7427603292077943,3,1,0.4291063333333334,1.287319,0.007971632676749163,0.007971632676749165,0.2933758520721111,2.640382668649,0.005921012511305819,0.00592101251130582),
    @("ECs","Lta","Tnfrsf14","Resolving-Mac",3,1,0.6836903333333334,2.051071,0.7427603292077943,0.7427603292077943,3,1,19.34413,58.03239,0.3593615074692841,0.3593615074692842,13.22539468774333,119.02855218969,0.2669194715924947,0.2669194715924948),
    @("Inflammatory-Mac","Lta","Tnfrsf14","ECs",1,0.3333333333333333,0.174178,0.522534,0.1892267629264251,0.1892267629264251,3,1,5.578493666666667,16.735481,0.1036332930693284,0.1036332930693284,0.9716508698726668,8.744857828854,0.01961019257891454,0.01961019257891454),
    @("Inflammatory-Mac","Lta","Tnfrsf14","FAPs",1,0.3333333333333333,0.174178,0.522534,0.1892267629264251,0.1892267629264251,3,1,3.704475666666667,11.113427,0.06881911773528272,0.06881911773528274,0.6452381626686667,5.807143464018001,0.01302241887650008,0.01302241887650008),
    @("Inflammatory-Mac","Lta","Tnfrsf14","Inflammatory-Mac",1,0.3333333333333333,0.174178,0.522534,0.1892267629264251,0.1892267629264251,3,1,24.77295966666667,74.31887900000001,0.4602144490493554,0.4602144490493556,4.314904568820667,38.83414111938601,0.08708489044557773,0.08708489044557774),
    @("Inflammatory-Mac","Lta","Tnfrsf14","MuSCs",1,0.3333333333333333,0.174178,0.522534,0.1892267629264251,0.1892267629264251,3,1,0.4291063333333334,1.287319,0.007971632676749163,0.007971632676749165,0.07474088292733336,0.6726679463460001,0.001508446246659757,0.001508446246659757),
    @("Inflammatory-Mac","Lta","Tnfrsf14","Resolving-Mac",1,0.3333333333333333,0.174178,0.522534,0.1892267629264251,0.1892267629264251,3,1,19.34413,58.03239,0.3593615074692841,0.3593615074692842,3.369321875140001,30.32389687626,0.06800081477877297,0.06800081477877297),
    @("Resolving-Mac","Lta","Tnfrsf14","ECs",2,0.6666666666666666,0.06260399999999999,0.187812,0.06801290786578049,0.06801290786578047,3,1,5.578493666666667,16.735481,0.1036332930693284,0.1036332930693284,0.349236017508,3.143124157572,0.00704840161335166,0.00704840161335166),
    @("Resolving-Mac","Lta","Tnfrsf14","FAPs",2,0.6666666666666666,0.06260399999999999,0.187812,0.06801290786578049,0.06801290786578047,3,1,3.704475666666667,11.113427,0.06881911773528272,0.06881911773528274,0.231914994636,2.087234951724,0.004680588313934084,0.004680588313934084),
    @("Resolving-Mac","Lta","Tnfrsf14","Inflammatory-Mac",2,0.6666666666666666,0.06260399999999999,0.187812,0.06801290786578049,0.06801290786578047,3,1,24.77295966666667,74.31887900000001,0.4602144490493554,0.4602144490493556,1.550886366972,13.957977302748,0.03130052292169474,0.03130052292169475),
    @("Resolving-Mac","Lta","Tnfrsf14","MuSCs",2,0.6666666666666666,0.06260399999999999,0.187812,0.06801290786578049,0.06801290786578047,3,1,0.4291063333333334,1.287319,0.007971632676749163,0.007971632676749165,0.026863772892,0.241773956028,0.0005421739187835858,0.000542173918783586),
    @("Resolving-Mac","Lta","Tnfrsf14","Resolving-Mac",2,0.6666666666666666,0.06260399999999999,0.187812,0.06801290786578049,0.06801290786578047,3,1,19.34413,58.03239,0.3593615074692841,0.3593615074692842,1.21101991452,10.89917923068,0.02444122109801641,0.02444122109801641)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

